# working on vaporize_theia.py mg/si and mg/al
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Table 1 (sheet1) - updated thermodynamic values for R1 (row 2), R9 (row 10),
# R10 (row 11), R11 (row 12), R12 (row 13); new column I ("C") added with a
# value for row 2.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Table 1")

$ws1.Range("I1").Value = "C"

$ws1.Range("G2").Value = 23.278899717000002
$ws1.Range("H2").Value = -99057.305063000007
$ws1.Range("I2").Value = 2655395.9

$ws1.Range("G10").Value = -9.5500000000000007
$ws1.Range("H10").Value = 63948

$ws1.Range("G11").Value = -26.91
$ws1.Range("H11").Value = 204359

$ws1.Range("G12").Value = -29.86
$ws1.Range("H12").Value = 200903

$ws1.Range("H13").Value = -33554.101999999999

$null = $ws1.Range("H14").Select()

# ---------------------------------------------------------------------------
# Table 3 (sheet3) - refined fit coefficients on rows 2 and 3.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Table 3")

$ws3.Range("C2").Value = -0.94840000000000002
$ws3.Range("D2").Value = 7404.5

$ws3.Range("C3").Value = 0.41520000000000001
$ws3.Range("D3").Value = 2330.6170000000002

$null = $ws3.Range("F2").Select()

# ---------------------------------------------------------------------------
# Table 4 (sheet4) - a new "C" column is inserted before the existing
# "Reactants" column (which shifts from E to F). Row 3 (R2 reaction) gets a
# value in the new column plus updated fit coefficients.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Table 4")

$ws4.Range("E:E").Insert()
$ws4.Range("E1").Value = "C"

$ws4.Range("C3").Value = -12.8489
$ws4.Range("D3").Value = 43367.3
$ws4.Range("E3").Formula = "=-2655400"

$null = $ws4.Range("J11").Select()

# Table 4 is the tab that was active/selected when the workbook was saved.
$null = $ws4.Activate()
